# Auto-generated market-price refresh for Golem_Profits workbook
# Updates currentAveragePrice*/Leve cost columns (H..N) per sheet/row as
# produced by the scheduled price-sync runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H17" = 8356.643
    "J17" = 10899.3
    "L17" = 32697.9
    "N17" = -33033.89999999999
    "H70" = 3358.3333
    "I70" = 3481.818
    "K70" = 10445.454
    "M70" = -10175.454
    "H73" = 3358.3333
    "I73" = 3481.818
    "K73" = 10445.454
    "M73" = -9509.454000000002
    "H131" = 1562
    "I131" = 1562
    "K131" = 4686
    "M131" = 354
    "H133" = 200000
    "I133" = 200000
    "K133" = 200000
    "M133" = -194940
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H45" = 2161.25
    "I45" = 2161.25
    "K45" = 2161.25
    "M45" = -1784.25
    "H102" = 4745.3335
    "I102" = 4745.3335
    "K102" = 4745.3335
    "M102" = -3123.3335
    "H122" = 2515.4614
    "I122" = 2600.0833
    "K122" = 7800.249899999999
    "M122" = -5350.249899999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H94" = 2450
    "J94" = 2900
    "L94" = 2900
    "N94" = -3802
    "H107" = 1455.3846
    "I107" = 1455.3846
    "K107" = 1455.3846
    "M107" = 464.6153999999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H5" = 4988.3335
    "I5" = 5978.5
    "K5" = 5978.5
    "M5" = -5866.5
    "H31" = 8249.5
    "I31" = 7332.778
    "K31" = 7332.778
    "M31" = -7037.778
    "H34" = 8249.5
    "I34" = 7332.778
    "K34" = 7332.778
    "M34" = -7130.778
    "H39" = 21399.8
    "I39" = 8999.666999999999
    "K39" = 8999.666999999999
    "M39" = -8608.666999999999
    "H49" = 21399.8
    "I49" = 8999.666999999999
    "K49" = 8999.666999999999
    "M49" = -8817.666999999999
    "H75" = 15260
    "J75" = 15260
    "L75" = 15260
    "N75" = -17256
    "H78" = 15260
    "J78" = 15260
    "L78" = 45780
    "N78" = -55764
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H23" = 79
    "J23" = 0
    "L23" = 0
    "H38" = 1178.5714
    "I38" = 500
    "J38" = 1687.5
    "K38" = 1500
    "L38" = 5062.5
    "M38" = -1153
    "N38" = -5756.5
    "H86" = 2386.875
    "I86" = 2270.7144
    "J86" = 3200
    "K86" = 6812.1432
    "L86" = 9600
    "M86" = -5626.1432
    "N86" = -11972
    "H89" = 2386.875
    "I89" = 2270.7144
    "J89" = 3200
    "K89" = 20436.4296
    "L89" = 28800
    "M89" = -14508.4296
    "N89" = -40656
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("N23")) {
    $ws.Range($addr).ClearContents()
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H7" = 2833333
    "H8" = 2833333
    "H11" = 4835770
    "I11" = 6740556.5
    "K11" = 6740556.5
    "M11" = -6740417.5
    "H14" = 5644557.5
    "I14" = 5644557.5
    "J14" = 0
    "K14" = 5644557.5
    "L14" = 0
    "M14" = -5644389.5
    "H122" = 6003.5
    "I122" = 6003.5
    "K122" = 18010.5
    "M122" = -15560.5
    "H132" = 1999
    "I132" = 1999
    "K132" = 5997
    "M132" = -3467
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("N14")) {
    $ws.Range($addr).ClearContents()
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H16" = 2500
    "I16" = 2500
    "K16" = 2500
    "M16" = -2330
    "H61" = 20741
    "I61" = 900
    "J61" = 50502.5
    "K61" = 900
    "L61" = 50502.5
    "M61" = -698
    "N61" = -50906.5
    "H68" = 1149.5
    "I68" = 1799
    "J68" = 500
    "K68" = 1799
    "L68" = 500
    "M68" = -1050
    "N68" = -1998
    "H69" = 0
    "J69" = 0
    "L69" = 0
    "H71" = 1149.5
    "I71" = 1799
    "J71" = 500
    "K71" = 8995
    "L71" = 2500
    "M71" = -5251
    "N71" = -9988
    "H72" = 0
    "J72" = 0
    "L72" = 0
    "H113" = 20741
    "I113" = 900
    "J113" = 50502.5
    "K113" = 900
    "L113" = 50502.5
    "M113" = 1270
    "N113" = -54842.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("N69", "N72")) {
    $ws.Range($addr).ClearContents()
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H120" = 565315.5
    "J120" = 565315.5
    "L120" = 565315.5
    "N120" = -574991.5
    "H121" = 0
    "J121" = 0
    "L121" = 0
    "H122" = 0
    "I122" = 0
    "J122" = 0
    "K122" = 0
    "L122" = 0
    "H124" = 0
    "J124" = 0
    "L124" = 0
    "H125" = 0
    "J125" = 0
    "L125" = 0
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("N121", "M122", "N122", "N124", "N125")) {
    $ws.Range($addr).ClearContents()
}
